# Code_Verification_Status.xlsx update
# Reflects: removal of "Files to be fixed" style highlighting on the
# "Files to be created" block, insertion of two newly-renamed script rows
# (3a_145879.sh, and a restructured order of the renamed-script list),
# and a text correction from "5_profiling_all.R" to "5_profiling.R".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XX_down_samp")

# --- Row 10/11: remove the amber "Files to be fixed" style, they become
#     plain "Files to be created" rows (col A default style, col B:E style
#     matching the bold-centered style used elsewhere, e.g. row 13). ---
$ws.Range("A10:A11").ClearFormats()

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B10:E10").PasteSpecial(-4122) | Out-Null
# B11:E11 already uses the correct style (s="3"); leave untouched.

# --- Insert one new row at position 15 so the renamed-script list grows
#     from 5 rows (14-18) to 6 rows (14-19). ---
$ws.Rows("15:15").Insert()

# Row 14: 3a_run_mSigHdp.R -> 2a_run_mSigHdp.R  (wrap style, orange fill on B:E)
$ws.Range("A14").Value2 = "3a_run_mSigHdp.R -> 2a_run_mSigHdp.R"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B14:E14").PasteSpecial(-4122) | Out-Null
$c = $ws.Range("A14").Characters(22, 2)
$c.Font.Color = 255

# Row 15 (new row): 3a_145879.sh -> 2a_145879.sh
$ws.Range("A15").Value2 = "3a_145879.sh -> 2a_145879.sh"
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B15:E15").PasteSpecial(-4122) | Out-Null
$c = $ws.Range("A15").Characters(17, 2)
$c.Font.Color = 255

# Row 16: 5_Summarize.R -> 3_Summarize.R (no special style - default)
$ws.Range("A16").Value2 = "5_Summarize.R -> 3_Summarize.R"
$ws.Range("A16").ClearFormats()
$c = $ws.Range("A16").Characters(18, 1)
$c.Font.Color = 255

# Row 17: 6_combine_extraction_metrics.R -> 4_combine_extraction_metrics.R
$ws.Range("A17").Value2 = "6_combine_extraction_metrics.R ->`n4_combine_extraction_metrics.R"
$c = $ws.Range("A17").Characters(33, 1)
$c.Font.Color = 255

# Row 18: 7b_profiling_all.R -> 5_profiling.R  (note: renamed file name changed
# from "5_profiling_all.R" to "5_profiling.R")
$ws.Range("A18").Value2 = "7b_profiling_all.R ->`n5_profiling.R"
$c = $ws.Range("A18").Characters(23, 1)
$c.Font.Color = 255

# Row 19: 8_plotting.R -> 6_plotting.R
$ws.Range("A19").Value2 = "8_plotting.R -> 6_plotting.R"
$c = $ws.Range("A19").Characters(16, 1)
$c.Font.Color = 255

# --- Update selection / active cell to mirror the authored session ---
$ws.Range("A18").Select() | Out-Null
